$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain their text formatting so values
# such as "1.00" or "10.80" are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "72.561.99"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3
$ws.Range("D3").Value = "3.957.31"
$ws.Range("E3").Value = "  -0.87%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").Value = "581.18"
$ws.Range("E5").Value = "  +6.87%  "

# Row 6
$ws.Range("D6").Value = "156.69"
$ws.Range("E6").Value = "  +3.91%  "

# Row 7
$ws.Range("D7").Value = "0.677"
$ws.Range("E7").Value = "  -3.55%  "

# Row 8
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "0.745"
$ws.Range("E9").Value = "  -0.31%  "

# Row 10
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  -2.00%  "

# Row 11
$ws.Range("D11").Value = "53.77"
$ws.Range("E11").Value = "  +0.88%  "

# Row 12
$ws.Range("D12").Value = "0.0000317"
$ws.Range("E12").Value = "  -2.08%  "

# Row 13
$ws.Range("D13").Value = "10.80"
$ws.Range("E13").Value = "  +1.04%  "

# Row 14
$ws.Range("D14").Value = "4.608.77"
$ws.Range("E14").Value = "  -0.59%  "

# Row 15
$ws.Range("D15").Value = "3.963.60"
$ws.Range("E15").Value = "  -1.17%  "

# Row 16
$ws.Range("D16").Value = "1.27"
$ws.Range("E16").Value = "  +7.18%  "

# Row 17
$ws.Range("D17").Value = "13.95"
$ws.Range("E17").Value = "  -1.50%  "

# Row 18
$ws.Range("D18").Value = "20.38"
$ws.Range("E18").Value = "  -0.98%  "

# Row 19
$ws.Range("E19").Value = "  -0.46%  "

# Row 20
$ws.Range("D20").Value = "72.460.03"
$ws.Range("E20").Value = "  +0.88%  "

# Row 21
$ws.Range("D21").Value = "431.10"
$ws.Range("E21").Value = "  -0.12%  "

# Row 22
$ws.Range("D22").Value = "4.66"
$ws.Range("E22").Value = "  +8.41%  "

# Row 23
$ws.Range("D23").Value = "95.66"
$ws.Range("E23").Value = "  -1.28%  "

# Row 24
$ws.Range("E24").Value = "  -3.70%  "

# Row 25
$ws.Range("D25").Value = "14.28"
$ws.Range("E25").Value = "  -0.42%  "

# Row 26
$ws.Range("D26").Value = "4.43"
$ws.Range("E26").Value = "  +22.35%  "

# Row 27
$ws.Range("D27").Value = "11.17"
$ws.Range("E27").Value = "  -3.13%  "

# Row 28
$ws.Range("D28").Value = "10.71"
$ws.Range("E28").Value = "  +0.02%  "

# Row 29
$ws.Range("E29").Value = "  +1.48%  "

# Row 30
$ws.Range("D30").Value = "36.25"
$ws.Range("E30").Value = "  -1.37%  "

# Row 31
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  +4.66%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "13.58"
$ws.Range("E32").Value = "  +0.75%  "

# Row 33
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "49.89"
$ws.Range("E33").Value = "  +1.83%  "

# Row 34
$ws.Range("E34").Value = "  -0.47%  "

# Row 35
$ws.Range("D35").Value = "677.39"
$ws.Range("E35").Value = "  -0.22%  "

# Row 36
$ws.Range("D36").Value = "68.59"
$ws.Range("E36").Value = "  +4.05%  "

# Row 37
$ws.Range("D37").Value = "0.435"
$ws.Range("E37").Value = "  -2.24%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0854"
$ws.Range("E38").Value = "  +2.30%  "

# Row 39
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  +0.04%  "

# Row 40
$ws.Range("E40").Value = "  -0.20%  "

# Row 41
$ws.Range("D41").Value = "0.145"
$ws.Range("E41").Value = "  -4.05%  "

# Row 42
$ws.Range("D42").Value = "3.32"
$ws.Range("E42").Value = "  -1.78%  "

# Row 43
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("D44").Value = "10.88"
$ws.Range("E44").Value = "  +10.38%  "

# Row 45
$ws.Range("D45").Value = "0.0485"
$ws.Range("E45").Value = "  -0.44%  "

# Row 46
$ws.Range("D46").Value = "0.148"
$ws.Range("E46").Value = "  -1.22%  "

# Row 47
$ws.Range("D47").Value = "2.67"
$ws.Range("E47").Value = "  -3.46%  "

# Row 48
$ws.Range("D48").Value = "3.40"
$ws.Range("E48").Value = "  +0.94%  "

# Row 49
$ws.Range("D49").Value = "3.41"
$ws.Range("E49").Value = "  +3.96%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "2.99"
$ws.Range("E50").Value = "  -0.07%  "

# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "2.15"
$ws.Range("E51").Value = "  +6.71%  "

